{"js": "// Proof-reading pass over the \"Justification of consolidated changes\" text:\n//  - tighten/clarify wording in several sentences\n//  - move the stray \"_GoBack\" bookmark from the end of the last paragraph\n//    into the \"fill in gaps,\" paragraph, right after the comma.\n\nasync function replaceOnce(findText, replaceText) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"get rid of\" -> \"remove\"\nawait replaceOnce(\n  \"get rid of the OFFICE table\",\n  \"remove the OFFICE table\"\n);\n\n// 2) Tidy up the OFFICE/STAFF/LOCATION sentence pair.\nawait replaceOnce(\n  \"duplicated in the OFFICE table that he had. Another change that was related to the STAFF and OFFICE relationship was the relationship that all of us had between STAFF and LOCATION.\",\n  \"duplicated in the OFFICE table. Another change was to relate the STAFF and OFFICE relationship to a relationship between STAFF and LOCATION consistent with the rest of the group\\u2019s feedback.\"\n);\n\n// 3) \"Another major change that we made when making the final ERD was to\n//    accommodate\" -> \"Changes were made to accommodate and relate\"\nawait replaceOnce(\n  \"Another major change that we made when making the final ERD was to accommodate both events\",\n  \"Changes were made to accommodate and relate both events\"\n);\n\n// 4) \"integrated them together.\" -> \"integrated the best features.\"\nawait replaceOnce(\"integrated them together.\", \"integrated the best features.\");\n\n// 5) \"fill in gaps such as some of the missing attributes\"\n//    -> \"fill in gaps, missing attributes\"\nawait replaceOnce(\n  \"fill in gaps such as some of the missing attributes\",\n  \"fill in gaps, missing attributes\"\n);\n\n// 6) Relocate the \"_GoBack\" bookmark: delete it from the end of the last\n//    real paragraph and re-insert it right after \"fill in gaps,\".\nconst bookmarkName = \"_GoBack\";\nconst existing = context.document.getBookmarkRangeOrNullObject(bookmarkName);\nexisting.load(\"isNullObject\");\nawait context.sync();\nif (!existing.isNullObject) {\n  context.document.deleteBookmark(bookmarkName);\n  await context.sync();\n}\n\nconst anchorResults = context.document.body.search(\"fill in gaps,\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Anchor text for bookmark not found\");\n}\nconst anchorRange = anchorResults.items[0].getRange(\"End\");\nanchorRange.insertBookmark(bookmarkName);\nawait context.sync();\n", "ps1": "# Proof-reading pass over the \"Justification of consolidated changes\" text:\n#  - tighten/clarify wording in several sentences\n#  - move the stray \"_GoBack\" bookmark from the end of the last paragraph\n#    into the \"fill in gaps,\" paragraph, right after the comma.\n\n$doc = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $doc.Content\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1) \"get rid of\" -> \"remove\"\nReplace-Text \"get rid of the OFFICE table\" \"remove the OFFICE table\"\n\n# 2) Tidy up the OFFICE/STAFF/LOCATION sentence pair.\nReplace-Text \"duplicated in the OFFICE table that he had. Another change that was related to the STAFF and OFFICE relationship was the relationship that all of us had between STAFF and LOCATION.\" \"duplicated in the OFFICE table. Another change was to relate the STAFF and OFFICE relationship to a relationship between STAFF and LOCATION consistent with the rest of the group\u2019s feedback.\"\n\n# 3) \"Another major change that we made when making the final ERD was to\n#    accommodate\" -> \"Changes were made to accommodate and relate\"\nReplace-Text \"Another major change that we made when making the final ERD was to accommodate both events\" \"Changes were made to accommodate and relate both events\"\n\n# 4) \"integrated them together.\" -> \"integrated the best features.\"\nReplace-Text \"integrated them together.\" \"integrated the best features.\"\n\n# 5) \"fill in gaps such as some of the missing attributes\"\n#    -> \"fill in gaps, missing attributes\"\nReplace-Text \"fill in gaps such as some of the missing attributes\" \"fill in gaps, missing attributes\"\n\n# 6) Relocate the \"_GoBack\" bookmark: delete it from the end of the last\n#    real paragraph and re-insert it right after \"fill in gaps,\".\nif ($doc.Bookmarks.Exists(\"_GoBack\")) {\n    $doc.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $doc.Content\n$found = $anchor.Find.Execute(\"fill in gaps,\")\nif (-not $found) {\n    throw \"Anchor text for bookmark not found\"\n}\n$anchor.Collapse(0)\n$doc.Bookmarks.Add(\"_GoBack\", $anchor)\n"}
